# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# per the scraped-data refresh (GitHub Actions cron run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.851.69"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.635.88"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.861.85"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "1.634.62"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5602"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "0.0₅7667"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "25.856.71"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.388"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.950"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.150"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.790"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1231"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.851"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04959"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.300"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.566"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9035"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.577"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").Value = "1.134.28"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01569"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9965"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8011"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "0.0₈112"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4264"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.782"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05056"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +0.18%  "
